# Add "Wins", "Losses", "Ties" columns (AC, AD, AE) holding the team's
# season record, mirroring the header formatting already used by the
# other header cells (e.g. column AB).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header cell format onto the three new header cells.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Populate every player/data row (2-48) with the season record.
for ($row = 2; $row -le 48; $row++) {
    $ws.Cells.Item($row, 29).Value = 61   # AC: Wins
    $ws.Cells.Item($row, 30).Value = 101  # AD: Losses
    $ws.Cells.Item($row, 31).Value = 0    # AE: Ties
}
